# Applies the XML diff to the document: adds " Alex!!!" + moved bookmark to
# the Heading1 paragraph, splits/rewrites several runs, inserts extra
# bulleted items, removes the stray proofErr-wrapped "statement" run
# (replacing it with "statements."), and appends a run of new trailing
# paragraphs.
#
# Strategy: the whole body content (everything the old `$d.Content` range
# covers, i.e. every paragraph up to - but not including - the final
# `sectPr`) is replaced in one shot via `Range.InsertXML` using a
# `pkg:package` WordProcessingML fragment that expresses the complete
# target body. This guarantees exact run-splitting (several diff hunks
# split a single run into multiple sibling `<w:r>` elements) which plain
# Find/Replace cannot produce, while `InsertXML` on `$d.Content` still
# leaves the trailing `<w:sectPr>` (page setup) untouched.

$d = $word.ActiveDocument

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Translation items</w:t></w:r><w:r><w:t xml:space="preserve"> Alex!!!</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:t xml:space="preserve">A first segment for translation. </w:t></w:r><w:r><w:t>A second segment for translation</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>And another segment for translation</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Bulleted item </w:t></w:r><w:r><w:t>one</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Bulleted item </w:t></w:r><w:r><w:t>two</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Bulleted item </w:t></w:r><w:r><w:t>one</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Bulleted item three</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Bulleted item three</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">My closing </w:t></w:r><w:r><w:t>statements.</w:t></w:r></w:p><w:p><w:r><w:t>Let’s try another repetition. And another segment for translation</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>A first segment for translation.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Another sentence here</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Text 5</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Translation provider.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>This is a new segment.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Content.InsertXML($xml)
